$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.138.42"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "'1.858.24"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'234.29"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").Value = "'0.2822"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "'0.06555"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'20.20"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("D11").Value = "'0.07802"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'97.32"
$ws.Range("E12").Value = "  -5.94%  "
$ws.Range("D13").Value = "'1.859.89"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "'5.093"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'285.29"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").Value = "'30.171.73"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'5.459"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'2.098.80"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").Value = "'0.000007250"
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'6.147"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "'167.98"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").Value = "'9.315"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'19.08"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'1.931"
$ws.Range("E28").Value = "  -7.26%  "
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").Value = "'0.09644"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "'4.414"
$ws.Range("E31").Value = "  -3.86%  "
$ws.Range("D32").Value = "'1.471"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").Value = "'4.110"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").Value = "'0.04685"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "'0.6982"
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("D36").Value = "'1.088"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "'0.9995"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'2.704"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'0.01863"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").Value = "'6.309"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").Value = "'2.509"
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("D42").Value = "'72.12"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").Value = "'0.8621"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'1.944"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "'104.36"
$ws.Range("D46").Value = "'0.4165"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").Value = "'0.9999"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'1.028.72"
$ws.Range("E48").Value = "  +7.85%  "
$ws.Range("D49").Value = "'7.243"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").Value = "'9.157"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("D51").Value = "'33.83"
$ws.Range("E51").Value = "  -2.42%  "
